$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $v = $cell.Value2
    if ($v -eq $null) {
        continue
    }
    $parts = $v.Split(",")
    if ($parts.Length -lt 2) {
        continue
    }
    $trimmed = @()
    foreach ($p in $parts) {
        $trimmed += $p.Trim()
    }
    $last = $trimmed[$trimmed.Length - 1]
    if ($last -eq "System") {
        $rest = $trimmed[0..($trimmed.Length - 2)]
        $newParts = @("System") + $rest
        $newValue = [string]::Join(", ", $newParts)
        $cell.Value = $newValue
    }
}
